# Auto-generated Excel COM-interop script to apply scheduled market-data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1307407.8
$ws.Range("J17").Value = 1307407.8
$ws.Range("L17").Value = 3922223.4
$ws.Range("N17").Value = -3922559.4
$ws.Range("H28").Value = 655.9286
$ws.Range("I28").Value = 562.36365
$ws.Range("K28").Value = 562.36365
$ws.Range("M28").Value = -77.36365000000001
$ws.Range("H41").Value = 323.375
$ws.Range("I41").Value = 196.33333
$ws.Range("K41").Value = 196.33333
$ws.Range("M41").Value = 243.66667
$ws.Range("H74").Value = 4379.7
$ws.Range("I74").Value = 3599.625
$ws.Range("J74").Value = 7500
$ws.Range("K74").Value = 3599.625
$ws.Range("L74").Value = 7500
$ws.Range("M74").Value = -2663.625
$ws.Range("N74").Value = -9372
$ws.Range("H77").Value = 4379.7
$ws.Range("I77").Value = 3599.625
$ws.Range("J77").Value = 7500
$ws.Range("K77").Value = 17998.125
$ws.Range("L77").Value = 37500
$ws.Range("M77").Value = -13318.125
$ws.Range("N77").Value = -46860
$ws.Range("H116").Value = 6914.5557
$ws.Range("I116").Value = 4850
$ws.Range("J116").Value = 8566.200000000001
$ws.Range("K116").Value = 4850
$ws.Range("L116").Value = 8566.200000000001
$ws.Range("M116").Value = -1408
$ws.Range("N116").Value = -15450.2
$ws.Range("H137").Value = 12034.737
$ws.Range("I137").Value = 1297.6364
$ws.Range("K137").Value = 3892.9092
$ws.Range("M137").Value = -1342.9092
$ws.Range("H138").Value = 2552.074
$ws.Range("I138").Value = 3673.6875
$ws.Range("J138").Value = 2079.8157
$ws.Range("K138").Value = 11021.0625
$ws.Range("L138").Value = 6239.4471
$ws.Range("M138").Value = -5881.0625
$ws.Range("N138").Value = -16519.4471

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H47").Value = 13750
$ws.Range("J47").Value = 17333.334
$ws.Range("L47").Value = 17333.334
$ws.Range("N47").Value = -18783.334
$ws.Range("H61").Value = 84188.11
$ws.Range("I61").Value = 2385.3
$ws.Range("J61").Value = 288695.12
$ws.Range("K61").Value = 2385.3
$ws.Range("L61").Value = 288695.12
$ws.Range("M61").Value = -2173.3
$ws.Range("N61").Value = -289119.12
$ws.Range("H63").Value = 2579.5
$ws.Range("J63").Value = 2277.6667
$ws.Range("L63").Value = 2277.6667
$ws.Range("N63").Value = -3649.6667
$ws.Range("H66").Value = 2579.5
$ws.Range("J66").Value = 2277.6667
$ws.Range("L66").Value = 11388.3335
$ws.Range("N66").Value = -18252.3335
$ws.Range("H122").Value = 2401.7144
$ws.Range("I122").Value = 2317.2307
$ws.Range("K122").Value = 6951.6921
$ws.Range("M122").Value = -4501.6921
$ws.Range("H132").Value = 2048240.2
$ws.Range("I132").Value = 1168.5122
$ws.Range("J132").Value = 12539483
$ws.Range("K132").Value = 3505.536599999999
$ws.Range("L132").Value = 37618449
$ws.Range("M132").Value = -975.5365999999995
$ws.Range("N132").Value = -37623509
$ws.Range("H136").Value = 84188.11
$ws.Range("I136").Value = 2385.3
$ws.Range("J136").Value = 288695.12
$ws.Range("K136").Value = 7155.900000000001
$ws.Range("L136").Value = 866085.36
$ws.Range("M136").Value = -4605.900000000001
$ws.Range("N136").Value = -871185.36

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1985.3334
$ws.Range("I86").Value = 1494.875
$ws.Range("J86").Value = 2377.7
$ws.Range("K86").Value = 1494.875
$ws.Range("L86").Value = 2377.7
$ws.Range("M86").Value = -371.875
$ws.Range("N86").Value = -4623.7
$ws.Range("H89").Value = 1985.3334
$ws.Range("I89").Value = 1494.875
$ws.Range("J89").Value = 2377.7
$ws.Range("K89").Value = 7474.375
$ws.Range("L89").Value = 11888.5
$ws.Range("M89").Value = -1858.375
$ws.Range("N89").Value = -23120.5
$ws.Range("H94").Value = 1631.9783
$ws.Range("I94").Value = 1331.7576
$ws.Range("J94").Value = 2394.077
$ws.Range("K94").Value = 1331.7576
$ws.Range("L94").Value = 2394.077
$ws.Range("M94").Value = -880.7575999999999
$ws.Range("N94").Value = -3296.077
$ws.Range("H99").Value = 13421.714
$ws.Range("I99").Value = 15503.353
$ws.Range("J99").Value = 4574.75
$ws.Range("K99").Value = 15503.353
$ws.Range("L99").Value = 4574.75
$ws.Range("M99").Value = -14005.353
$ws.Range("N99").Value = -7570.75
$ws.Range("H107").Value = 1923.5588
$ws.Range("I107").Value = 2113.087
$ws.Range("J107").Value = 1527.2727
$ws.Range("K107").Value = 2113.087
$ws.Range("L107").Value = 1527.2727
$ws.Range("M107").Value = -193.087
$ws.Range("N107").Value = -5367.2727
$ws.Range("H134").Value = 42818.805
$ws.Range("I134").Value = 51563.4
$ws.Range("K134").Value = 154690.2
$ws.Range("M134").Value = -152155.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1645.8
$ws.Range("I22").Value = 983.1667
$ws.Range("J22").Value = 2087.5557
$ws.Range("K22").Value = 983.1667
$ws.Range("L22").Value = 2087.5557
$ws.Range("M22").Value = -633.1667
$ws.Range("N22").Value = -2787.5557
$ws.Range("H68").Value = 40000
$ws.Range("J68").Value = 40000
$ws.Range("L68").Value = 40000
$ws.Range("N68").Value = -41498
$ws.Range("H71").Value = 40000
$ws.Range("J71").Value = 40000
$ws.Range("L71").Value = 120000
$ws.Range("N71").Value = -127488
$ws.Range("H122").Value = 2093.3333
$ws.Range("I122").Value = 1442.2778
$ws.Range("K122").Value = 4326.8334
$ws.Range("M122").Value = -1876.8334
$ws.Range("H132").Value = 23258132
$ws.Range("I132").Value = 2354.158
$ws.Range("J132").Value = 200002050
$ws.Range("K132").Value = 7062.474
$ws.Range("L132").Value = 600006150
$ws.Range("M132").Value = -4532.474
$ws.Range("N132").Value = -600011210
$ws.Range("H134").Value = 20412832
$ws.Range("I134").Value = 2139.7354
$ws.Range("K134").Value = 6419.206200000001
$ws.Range("M134").Value = -3884.206200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 181.20833
$ws.Range("I2").Value = 230.07143
$ws.Range("J2").Value = 112.8
$ws.Range("K2").Value = 1380.42858
$ws.Range("L2").Value = 676.8
$ws.Range("M2").Value = -1267.42858
$ws.Range("N2").Value = -902.8
$ws.Range("H6").Value = 63.333332
$ws.Range("I6").Value = 101
$ws.Range("J6").Value = 25.666666
$ws.Range("K6").Value = 303
$ws.Range("L6").Value = 76.99999800000001
$ws.Range("M6").Value = -190
$ws.Range("N6").Value = -302.999998
$ws.Range("H26").Value = 2594.125
$ws.Range("I26").Value = 2970.8
$ws.Range("J26").Value = 1966.3334
$ws.Range("K26").Value = 8912.400000000001
$ws.Range("L26").Value = 5899.0002
$ws.Range("M26").Value = -8624.400000000001
$ws.Range("N26").Value = -6475.0002
$ws.Range("H34").Value = 2692.2727
$ws.Range("I34").Value = 1976.2858
$ws.Range("J34").Value = 3945.25
$ws.Range("K34").Value = 5928.857400000001
$ws.Range("L34").Value = 11835.75
$ws.Range("M34").Value = -5844.857400000001
$ws.Range("N34").Value = -12003.75
$ws.Range("H39").Value = 7231
$ws.Range("J39").Value = 9539.799999999999
$ws.Range("L39").Value = 28619.4
$ws.Range("N39").Value = -29207.4
$ws.Range("H55").Value = 2659.4
$ws.Range("J55").Value = 2449
$ws.Range("L55").Value = 7347
$ws.Range("N55").Value = -7701
$ws.Range("H120").Value = 5165.5
$ws.Range("I120").Value = 5165.5
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 15496.5
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -10658.5
$ws.Range("N120").ClearContents()
$ws.Range("H131").Value = 1455.21
$ws.Range("J131").Value = 1477.8229
$ws.Range("L131").Value = 4433.468699999999
$ws.Range("N131").Value = -14513.4687

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5047.25
$ws.Range("I70").Value = 4112.25
$ws.Range("J70").Value = 5358.9165
$ws.Range("K70").Value = 4112.25
$ws.Range("L70").Value = 5358.9165
$ws.Range("M70").Value = -3842.25
$ws.Range("N70").Value = -5898.9165
$ws.Range("H73").Value = 5047.25
$ws.Range("I73").Value = 4112.25
$ws.Range("J73").Value = 5358.9165
$ws.Range("K73").Value = 4112.25
$ws.Range("L73").Value = 5358.9165
$ws.Range("M73").Value = -3176.25
$ws.Range("N73").Value = -7230.9165
$ws.Range("H80").Value = 11280.823
$ws.Range("J80").Value = 15567.5
$ws.Range("L80").Value = 15567.5
$ws.Range("N80").Value = -17563.5
$ws.Range("H83").Value = 11280.823
$ws.Range("J83").Value = 15567.5
$ws.Range("L83").Value = 77837.5
$ws.Range("N83").Value = -87821.5
$ws.Range("H105").Value = 30109.8
$ws.Range("J105").Value = 30120
$ws.Range("L105").Value = 30120
$ws.Range("N105").Value = -37108
$ws.Range("H113").Value = 2077.5
$ws.Range("I113").Value = 1613.8889
$ws.Range("J113").Value = 6250
$ws.Range("K113").Value = 1613.8889
$ws.Range("L113").Value = 6250
$ws.Range("M113").Value = 556.1111000000001
$ws.Range("N113").Value = -10590
$ws.Range("H122").Value = 6386
$ws.Range("I122").Value = 7223.5557
$ws.Range("J122").Value = 4501.5
$ws.Range("K122").Value = 21670.6671
$ws.Range("L122").Value = 13504.5
$ws.Range("M122").Value = -19220.6671
$ws.Range("N122").Value = -18404.5
$ws.Range("H126").Value = 15992.75
$ws.Range("J126").Value = 3140
$ws.Range("L126").Value = 9420
$ws.Range("N126").Value = -14360
$ws.Range("I132").Value = 3779.5293
$ws.Range("J132").Value = 19168.5
$ws.Range("K132").Value = 11338.5879
$ws.Range("L132").Value = 57505.5
$ws.Range("M132").Value = -8808.5879
$ws.Range("N132").Value = -62565.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1772.5
$ws.Range("I40").Value = 1772.5
$ws.Range("K40").Value = 1772.5
$ws.Range("M40").Value = -1636.5
$ws.Range("H46").Value = 2068.125
$ws.Range("I46").Value = 498
$ws.Range("J46").Value = 2172.8
$ws.Range("K46").Value = 498
$ws.Range("L46").Value = 2172.8
$ws.Range("M46").Value = -310
$ws.Range("N46").Value = -2548.8
$ws.Range("H132").Value = 1091888.1
$ws.Range("I132").Value = 2471.9565
$ws.Range("J132").Value = 2881643.2
$ws.Range("K132").Value = 7415.869499999999
$ws.Range("L132").Value = 8644929.600000001
$ws.Range("M132").Value = -4885.869499999999
$ws.Range("N132").Value = -8649989.600000001
$ws.Range("H136").Value = 103245.305
$ws.Range("I136").Value = 15155.353
$ws.Range("K136").Value = 45466.05899999999
$ws.Range("M136").Value = -42916.05899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 421.10526
$ws.Range("I4").Value = 165.25
$ws.Range("J4").Value = 859.7143
$ws.Range("K4").Value = 165.25
$ws.Range("L4").Value = 859.7143
$ws.Range("M4").Value = -52.25
$ws.Range("N4").Value = -1085.7143
$ws.Range("H62").Value = 18366.889
$ws.Range("I62").Value = 13940.4
$ws.Range("J62").Value = 23900
$ws.Range("K62").Value = 13940.4
$ws.Range("L62").Value = 23900
$ws.Range("M62").Value = -13316.4
$ws.Range("N62").Value = -25148
$ws.Range("H65").Value = 18366.889
$ws.Range("I65").Value = 13940.4
$ws.Range("J65").Value = 23900
$ws.Range("K65").Value = 69702
$ws.Range("L65").Value = 119500
$ws.Range("M65").Value = -66582
$ws.Range("N65").Value = -125740
$ws.Range("H82").Value = 100300.5
$ws.Range("J82").Value = 100300.5
$ws.Range("L82").Value = 100300.5
$ws.Range("N82").Value = -101066.5
$ws.Range("H85").Value = 100300.5
$ws.Range("J85").Value = 100300.5
$ws.Range("L85").Value = 100300.5
$ws.Range("N85").Value = -102952.5
$ws.Range("H88").Value = 75125.664
$ws.Range("J88").Value = 80188.5
$ws.Range("L88").Value = 80188.5
$ws.Range("N88").Value = -81000.5
$ws.Range("H91").Value = 75125.664
$ws.Range("J91").Value = 80188.5
$ws.Range("L91").Value = 80188.5
$ws.Range("N91").Value = -82996.5
$ws.Range("H122").Value = 1920.0344
$ws.Range("I122").Value = 1899.3214
$ws.Range("K122").Value = 5697.9642
$ws.Range("M122").Value = -3247.9642
$ws.Range("H132").Value = 4882.2256
$ws.Range("I132").Value = 1764.8148
$ws.Range("K132").Value = 5294.4444
$ws.Range("M132").Value = -2764.4444
